# Apply the "456a3b4" gh-pages data refresh to both the "展览" sheet
# (Worksheets(1)) and the "全部类型" sheet (Worksheets(4)), which carry the
# same event list (the latter has one extra concert row inserted above
# row 32, which is why its row numbers diverge from row 32 onward).

$wb = $excel.ActiveWorkbook

# Row -> new "想去人数" (column F) value, identical deltas on both sheets
# except that every row from the old row 32 onward is shifted down by one
# row on the "全部类型" sheet.
$sheet1FUpdates = @{
    3  = 551
    4  = 1116
    5  = 104
    6  = 54
    9  = 1146
    10 = 16047
    11 = 258
    12 = 189
    13 = 1034
    14 = 6269
    16 = 121
    17 = 75
    21 = 6
    24 = 25
    25 = 20
    28 = 882
    29 = 35
    30 = 5026
    32 = 11209
    35 = 136
    36 = 195
    37 = 3825
    38 = 268
}

$sheet4FUpdates = @{
    3  = 551
    4  = 1116
    5  = 104
    6  = 54
    9  = 1146
    10 = 16047
    11 = 258
    12 = 189
    13 = 1034
    14 = 6269
    16 = 121
    17 = 75
    21 = 6
    24 = 25
    25 = 20
    28 = 882
    29 = 35
    30 = 5026
    33 = 11209
    36 = 136
    37 = 195
    38 = 3825
    39 = 268
}

function Update-EventSheet($ws, $fUpdates) {
    # Row 2 event ("苏州·国乙ony茶话会一对一委托-星渡咖啡") sold out / got
    # pulled: title now flagged "（取消）", 想去人数 ticks up by one, and
    # 最低票价 is replaced by the "不可售" (not sellable) status text.
    $ws.Range("C2").Value = "苏州·国乙ony茶话会一对一委托-星渡咖啡（取消）"
    $ws.Range("F2").Value = 3162
    $ws.Range("G2").Value = "不可售"

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }
}

Update-EventSheet $wb.Worksheets.Item(1) $sheet1FUpdates
Update-EventSheet $wb.Worksheets.Item(4) $sheet4FUpdates
